# Actualizacion automatica 2025-12-01 08:30:07
#
# 1) "VENTAS POR GRUPO": the current ("in progress") month's figures for a
#    number of asesor/cliente rows are reset to 0 now that the month has
#    closed out, and the "N de 21" completion counters in the totals row
#    (row 23) are recalculated to match.
# 2) "VENTA MENSUAL": the month columns roll forward by one (agosto drops
#    off, septiembre..diciembre take its place) so each month's figures
#    shift left by one column and a fresh (empty) month is appended on
#    the right; column widths follow the same shift.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO": zero-out the cells that rolled out of the
# "current" month bucket.
# ---------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$zeroCells = @(
    "D5", "L5", "M5",
    "D6", "E6", "L6",
    "H11", "I11", "M11",
    "M13",
    "D14", "I14", "K14", "L14", "M14",
    "H16", "M16",
    "M17",
    "D22", "H22", "L22", "M22"
)
foreach ($ref in $zeroCells) {
    $wsGrupo.Range($ref).Value = 0
}

# Row 23 "N de 21" completion counters.
$wsGrupo.Range("D23").Value = "0 de 21"
$wsGrupo.Range("E23").Value = "0 de 21"
$wsGrupo.Range("H23").Value = "0 de 21"
$wsGrupo.Range("I23").Value = "0 de 21"
$wsGrupo.Range("K23").Value = "0 de 21"
$wsGrupo.Range("L23").Value = "0 de 21"
$wsGrupo.Range("M23").Value = "0 de 21"

# ---------------------------------------------------------------------
# Sheet "VENTA MENSUAL": shift months (and column widths) one column to
# the left, dropping "agosto" and bringing in "diciembre" as a new,
# still-empty month.
# ---------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Column widths: C<-16, D<-14, E<-15 (F/G unchanged). ColumnWidth is in
# "characters", which Excel stores internally offset by +5/6 character;
# subtract that offset so the saved width matches exactly.
$widthOffset = 0.8333333333333333
$wsMensual.Columns.Item(3).ColumnWidth = 16 - $widthOffset
$wsMensual.Columns.Item(4).ColumnWidth = 14 - $widthOffset
$wsMensual.Columns.Item(5).ColumnWidth = 15 - $widthOffset

# Month header row.
$wsMensual.Range("C1").Value = "septiembre"
$wsMensual.Range("D1").Value = "octubre"
$wsMensual.Range("E1").Value = "noviembre"
$wsMensual.Range("F1").Value = "diciembre"

# Data rows: shift C<-D, D<-E, E<-F, F<-0 for every data row (5..23,
# skipping the blank rows that carry no month figures at all). Read the
# old values first (via the `.Value()` getter-call form) before writing
# anything, since C/D/E get overwritten in place.
$dataRows = 5..23
foreach ($r in $dataRows) {
    $cRange = $wsMensual.Range("C$r")
    $dRange = $wsMensual.Range("D$r")
    $eRange = $wsMensual.Range("E$r")
    $fRange = $wsMensual.Range("F$r")

    $newC = $dRange.Value()
    $newD = $eRange.Value()
    $newE = $fRange.Value()

    $cRange.Value = $newC
    $dRange.Value = $newD
    $eRange.Value = $newE
    $fRange.Value = 0
}
